# Apply cell updates for cryptos list refresh
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = "38.785.35"
$ws.Cells.Item(2, 5).Value = "  +0.31%  "

# Row 3
$ws.Cells.Item(3, 4).NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = "2.102.80"
$ws.Cells.Item(3, 5).Value = "  +0.06%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  +0.02%  "

# Row 5
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "228.62"
$ws.Cells.Item(5, 5).Value = "  -0.34%  "

# Row 6
$ws.Cells.Item(6, 5).Value = "  +0.14%  "

# Row 7
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "62.42"
$ws.Cells.Item(7, 5).Value = "  +1.43%  "

# Row 8
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.999"
$ws.Cells.Item(8, 5).Value = "  -0.03%  "

# Row 9
$ws.Cells.Item(9, 5).Value = "  +2.13%  "

# Row 10
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.0842"
$ws.Cells.Item(10, 5).Value = "  -0.31%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  -0.96%  "

# Row 12
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "15.74"
$ws.Cells.Item(12, 5).Value = "  +6.10%  "

# Row 13
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "2.413.43"
$ws.Cells.Item(13, 5).Value = "  +0.21%  "

# Row 14
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "22.10"
$ws.Cells.Item(14, 5).Value = "  -1.79%  "

# Row 15
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.812"
$ws.Cells.Item(15, 5).Value = "  +3.53%  "

# Row 16
$ws.Cells.Item(16, 5).Value = "  +0.44%  "

# Row 17
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "2.103.72"
$ws.Cells.Item(17, 5).Value = "  +0.40%  "

# Row 18
$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "38.820.08"
$ws.Cells.Item(18, 5).Value = "  +0.66%  "

# Row 19
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "71.96"
$ws.Cells.Item(19, 5).Value = "  +1.20%  "

# Row 20
$ws.Cells.Item(20, 5).Value = "  +0.88%  "

# Row 21
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "0.0₃0841"
$ws.Cells.Item(21, 5).Value = "  +0.53%  "

# Row 22
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "229.17"
$ws.Cells.Item(22, 5).Value = "  +0.84%  "

# Row 23
$ws.Cells.Item(23, 5).Value = "  +0.01%  "

# Row 24
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "2.36"
$ws.Cells.Item(24, 5).Value = "  -3.23%  "

# Row 25
$ws.Cells.Item(25, 5).Value = "  -0.19%  "

# Row 26
$ws.Cells.Item(26, 2).Value = "Cosmos"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "9.61"
$ws.Cells.Item(26, 5).Value = "  +1.76%  "

# Row 27
$ws.Cells.Item(27, 2).Value = "Monero"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "171.79"
$ws.Cells.Item(27, 5).Value = "  +0.85%  "

# Row 28
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "0.138"
$ws.Cells.Item(28, 5).Value = "  +5.37%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  +4.42%  "

# Row 30
$ws.Cells.Item(30, 5).Value = "  +1.10%  "

# Row 31
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "2.52"
$ws.Cells.Item(31, 5).Value = "  +8.75%  "

# Row 32
$ws.Cells.Item(32, 5).Value = "  +0.53%  "

# Row 33
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "4.56"
$ws.Cells.Item(33, 5).Value = "  +1.79%  "

# Row 34
$ws.Cells.Item(34, 5).Value = "  -0.56%  "

# Row 35
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "6.99"
$ws.Cells.Item(35, 5).Value = "  +7.19%  "

# Row 36
$ws.Cells.Item(36, 5).Value = "  +1.99%  "

# Row 37
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "2.43"
$ws.Cells.Item(37, 5).Value = "  +0.70%  "

# Row 38
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "3.58"
$ws.Cells.Item(38, 5).Value = "  +0.23%  "

# Row 39
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.00"
$ws.Cells.Item(39, 5).Value = "  -0.14%  "

# Row 40
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "18.06"
$ws.Cells.Item(40, 5).Value = "  -3.80%  "

# Row 41
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "102.94"
$ws.Cells.Item(41, 5).Value = "  +2.76%  "

# Row 42
$ws.Cells.Item(42, 5).Value = "  +3.47%  "

# Row 43
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "1.531.88"
$ws.Cells.Item(43, 5).Value = "  -1.01%  "

# Row 44
$ws.Cells.Item(44, 2).Value = "TrustWalletToken"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "1.17"
$ws.Cells.Item(44, 5).Value = "  +4.59%  "

# Row 45
$ws.Cells.Item(45, 2).Value = "FraxShare"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "7.87"
$ws.Cells.Item(45, 5).Value = "  +3.59%  "

# Row 46
$ws.Cells.Item(46, 5).Value = "  -1.19%  "

# Row 47
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.0912"
$ws.Cells.Item(47, 5).Value = "  -0.49%  "

# Row 48
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "4.12"
$ws.Cells.Item(48, 5).Value = "  -1.84%  "

# Row 49
$ws.Cells.Item(49, 5).Value = "  +1.10%  "

# Row 50
$ws.Cells.Item(50, 5).Value = "  -0.23%  "

# Row 51
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "2.300.15"
